$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("woman")

# --- Remove stray data that had been accidentally entered in row 37 ---
# (leave B37/C37 intact, clear D37/F37/H37/I37/J37)
$ws.Range('D37').ClearContents()
$ws.Range('F37').ClearContents()
$ws.Range('H37').ClearContents()
$ws.Range('I37').ClearContents()
$ws.Range('J37').ClearContents()

# --- Append new scraped rows (p_id 1-56) for the new account "soom._.inc" ---
$ws.Range('B103').Value = '1'
$ws.Range('C103').Value = 'soom._.inc'
$ws.Range('D103').Value = '1'
$ws.Range('F103').Value = '반팔'
$ws.Range('H103').Value = '면바지'
$ws.Range('I103').Value = '컨버스'
$ws.Range('J103').Value = '28'
$ws.Range('B104').Value = '2'
$ws.Range('C104').Value = 'soom._.inc'
$ws.Range('D104').Value = '7'
$ws.Range('F104').Value = '반팔'
$ws.Range('G104').Value = ' '
$ws.Range('H104').Value = '면바지'
$ws.Range('I104').Value = '단화'
$ws.Range('J104').Value = '28'
$ws.Range('B105').Value = '3'
$ws.Range('C105').Value = 'soom._.inc'
$ws.Range('D105').Value = '1'
$ws.Range('F105').Value = '반팔'
$ws.Range('H105').Value = '면바지'
$ws.Range('I105').Value = '단화'
$ws.Range('J105').Value = '28'
$ws.Range('B106').Value = '4'
$ws.Range('C106').Value = 'soom._.inc'
$ws.Range('D106').Value = '1,5'
$ws.Range('F106').Value = '나시'
$ws.Range('H106').Value = '면바지'
$ws.Range('I106').Value = '단화'
$ws.Range('J106').Value = '28'
$ws.Range('B107').Value = '5'
$ws.Range('C107').Value = 'soom._.inc'
$ws.Range('D107').Value = '4'
$ws.Range('F107').Value = '셔츠'
$ws.Range('H107').Value = '면바지'
$ws.Range('I107').Value = '구두'
$ws.Range('J107').Value = '19'
$ws.Range('B108').Value = '6'
$ws.Range('C108').Value = 'soom._.inc'
$ws.Range('D108').Value = '1'
$ws.Range('F108').Value = '맨투맨'
$ws.Range('H108').Value = '반바지'
$ws.Range('I108').Value = '운동화'
$ws.Range('J108').Value = '19'
$ws.Range('B109').Value = '7'
$ws.Range('C109').Value = 'soom._.inc'
$ws.Range('D109').Value = '1,5'
$ws.Range('F109').Value = '니트'
$ws.Range('H109').Value = '면바지'
$ws.Range('I109').Value = '운동화'
$ws.Range('J109').Value = '22'
$ws.Range('B110').Value = '8'
$ws.Range('C110').Value = 'soom._.inc'
$ws.Range('D110').Value = '1,5'
$ws.Range('F110').Value = '맨투맨'
$ws.Range('H110').Value = '면바지'
$ws.Range('I110').Value = '운동화'
$ws.Range('J110').Value = '19'
$ws.Range('B111').Value = '9'
$ws.Range('C111').Value = 'soom._.inc'
$ws.Range('D111').Value = '1,5'
$ws.Range('F111').Value = '조끼'
$ws.Range('H111').Value = '면바지'
$ws.Range('I111').Value = '컨버스'
$ws.Range('J111').Value = '28'
$ws.Range('B112').Value = '10'
$ws.Range('C112').Value = 'soom._.inc'
$ws.Range('D112').Value = '4,7'
$ws.Range('E112').Value = '블레이저'
$ws.Range('F112').Value = '반팔'
$ws.Range('H112').Value = '면바지'
$ws.Range('I112').Value = '구두'
$ws.Range('J112').Value = '12'
$ws.Range('B113').Value = '11'
$ws.Range('C113').Value = 'soom._.inc'
$ws.Range('D113').Value = '4,7'
$ws.Range('E113').Value = '자켓'
$ws.Range('H113').Value = '반바지'
$ws.Range('I113').Value = '컨버스'
$ws.Range('J113').Value = '12'
$ws.Range('B114').Value = '12'
$ws.Range('C114').Value = 'soom._.inc'
$ws.Range('D114').Value = '1,5,7'
$ws.Range('F114').Value = '니트'
$ws.Range('H114').Value = '면바지'
$ws.Range('I114').Value = '운동화'
$ws.Range('J114').Value = '19'
$ws.Range('B115').Value = '13'
$ws.Range('C115').Value = 'soom._.inc'
$ws.Range('D115').Value = '1,7'
$ws.Range('F115').Value = '나시'
$ws.Range('H115').Value = '면바지'
$ws.Range('I115').Value = '구두'
$ws.Range('J115').Value = '28'
$ws.Range('B116').Value = '14'
$ws.Range('C116').Value = 'soom._.inc'
$ws.Range('D116').Value = '7'
$ws.Range('F116').Value = '니트'
$ws.Range('G116').Value = '나시'
$ws.Range('H116').Value = '면바지'
$ws.Range('I116').Value = '운동화'
$ws.Range('J116').Value = '28'
$ws.Range('B117').Value = '15'
$ws.Range('C117').Value = 'soom._.inc'
$ws.Range('D117').Value = '1,5,7'
$ws.Range('F117').Value = '긴팔'
$ws.Range('H117').Value = '면바지'
$ws.Range('I117').Value = '운동화'
$ws.Range('J117').Value = '19'
$ws.Range('B118').Value = '16'
$ws.Range('C118').Value = 'soom._.inc'
$ws.Range('D118').Value = '5'
$ws.Range('E118').Value = '코트'
$ws.Range('H118').Value = '면바지'
$ws.Range('I118').Value = '운동화'
$ws.Range('J118').Value = '10'
$ws.Range('B119').Value = '17'
$ws.Range('C119').Value = 'soom._.inc'
$ws.Range('D119').Value = '1,5'
$ws.Range('F119').Value = '니트'
$ws.Range('H119').Value = '면바지'
$ws.Range('I119').Value = '컨버스'
$ws.Range('J119').Value = '19'
$ws.Range('B120').Value = '18'
$ws.Range('C120').Value = 'soom._.inc'
$ws.Range('D120').Value = '1,5'
$ws.Range('F120').Value = '니트'
$ws.Range('H120').Value = '면바지'
$ws.Range('I120').Value = '운동화'
$ws.Range('J120').Value = '19'
$ws.Range('B121').Value = '19'
$ws.Range('C121').Value = 'soom._.inc'
$ws.Range('D121').Value = '1,5,7'
$ws.Range('F121').Value = '긴팔'
$ws.Range('H121').Value = '면바지'
$ws.Range('I121').Value = '컨버스'
$ws.Range('J121').Value = '19'
$ws.Range('B122').Value = '20'
$ws.Range('C122').Value = 'soom._.inc'
$ws.Range('D122').Value = '1'
$ws.Range('F122').Value = '긴팔'
$ws.Range('H122').Value = '면바지'
$ws.Range('I122').Value = '구두'
$ws.Range('J122').Value = '19'
$ws.Range('B123').Value = '21'
$ws.Range('C123').Value = 'soom._.inc'
$ws.Range('D123').Value = '4'
$ws.Range('F123').Value = '니트'
$ws.Range('H123').Value = '면바지'
$ws.Range('I123').Value = '구두'
$ws.Range('J123').Value = '19'
$ws.Range('B124').Value = '22'
$ws.Range('C124').Value = 'soom._.inc'
$ws.Range('D124').Value = '1'
$ws.Range('F124').Value = '긴팔'
$ws.Range('H124').Value = '면바지'
$ws.Range('I124').Value = '운동화'
$ws.Range('J124').Value = '19'
$ws.Range('B125').Value = '23'
$ws.Range('C125').Value = 'soom._.inc'
$ws.Range('D125').Value = '1,5'
$ws.Range('F125').Value = '니트'
$ws.Range('H125').Value = '면바지'
$ws.Range('I125').Value = '컨버스'
$ws.Range('J125').Value = '19'
$ws.Range('B126').Value = '24'
$ws.Range('C126').Value = 'soom._.inc'
$ws.Range('D126').Value = '1'
$ws.Range('F126').Value = '긴팔'
$ws.Range('H126').Value = '면바지'
$ws.Range('I126').Value = '운동화'
$ws.Range('J126').Value = '19'
$ws.Range('B127').Value = '25'
$ws.Range('C127').Value = 'soom._.inc'
$ws.Range('D127').Value = '1,7'
$ws.Range('E127').Value = '블레이저'
$ws.Range('F127').Value = '긴팔'
$ws.Range('H127').Value = '면바지'
$ws.Range('I127').Value = '컨버스'
$ws.Range('J127').Value = '12'
$ws.Range('B128').Value = '26'
$ws.Range('C128').Value = 'soom._.inc'
$ws.Range('D128').Value = '1,7'
$ws.Range('E128').Value = '가디건'
$ws.Range('F128').Value = '긴팔'
$ws.Range('H128').Value = '면바지'
$ws.Range('I128').Value = '운동화'
$ws.Range('J128').Value = '19'
$ws.Range('B129').Value = '27'
$ws.Range('C129').Value = 'soom._.inc'
$ws.Range('D129').Value = '1'
$ws.Range('F129').Value = '긴팔'
$ws.Range('G129').Value = '긴팔'
$ws.Range('H129').Value = '면바지'
$ws.Range('I129').Value = '운동화'
$ws.Range('J129').Value = '19'
$ws.Range('B130').Value = '28'
$ws.Range('C130').Value = 'soom._.inc'
$ws.Range('D130').Value = '1,7'
$ws.Range('E130').Value = '바람막이'
$ws.Range('F130').Value = '반팔'
$ws.Range('H130').Value = '면바지'
$ws.Range('I130').Value = '운동화'
$ws.Range('J130').Value = '19'
$ws.Range('B131').Value = '29'
$ws.Range('C131').Value = 'soom._.inc'
$ws.Range('D131').Value = '1'
$ws.Range('F131').Value = '긴팔'
$ws.Range('H131').Value = '면바지'
$ws.Range('I131').Value = '운동화'
$ws.Range('J131').Value = '19'
$ws.Range('B132').Value = '30'
$ws.Range('C132').Value = 'soom._.inc'
$ws.Range('D132').Value = '5,7'
$ws.Range('E132').Value = '니트집업'
$ws.Range('F132').Value = '긴팔'
$ws.Range('H132').Value = '면바지'
$ws.Range('I132').Value = '운동화'
$ws.Range('J132').Value = '19'
$ws.Range('B133').Value = '31'
$ws.Range('C133').Value = 'soom._.inc'
$ws.Range('D133').Value = '1,5,7'
$ws.Range('F133').Value = '긴팔'
$ws.Range('G133').Value = '긴팔'
$ws.Range('H133').Value = '면바지'
$ws.Range('I133').Value = '컨버스'
$ws.Range('J133').Value = '19'
$ws.Range('B134').Value = '32'
$ws.Range('C134').Value = 'soom._.inc'
$ws.Range('D134').Value = '4'
$ws.Range('F134').Value = '원피스'
$ws.Range('I134').Value = '구두'
$ws.Range('J134').Value = '22'
$ws.Range('B135').Value = '33'
$ws.Range('C135').Value = 'soom._.inc'
$ws.Range('D135').Value = '1,7'
$ws.Range('F135').Value = '긴팔'
$ws.Range('H135').Value = '면바지'
$ws.Range('I135').Value = '운동화'
$ws.Range('J135').Value = '22'
$ws.Range('B136').Value = '34'
$ws.Range('C136').Value = 'soom._.inc'
$ws.Range('D136').Value = '1'
$ws.Range('F136').Value = '셔츠'
$ws.Range('G136').Value = '반팔'
$ws.Range('H136').Value = '면바지'
$ws.Range('I136').Value = '컨버스'
$ws.Range('J136').Value = '22'
$ws.Range('B137').Value = '35'
$ws.Range('C137').Value = 'soom._.inc'
$ws.Range('D137').Value = '1'
$ws.Range('F137').Value = '반팔'
$ws.Range('H137').Value = '면바지'
$ws.Range('I137').Value = '운동화'
$ws.Range('J137').Value = '28'
$ws.Range('B138').Value = '36'
$ws.Range('C138').Value = 'soom._.inc'
$ws.Range('D138').Value = '1'
$ws.Range('F138').Value = '긴팔'
$ws.Range('H138').Value = '면바지'
$ws.Range('I138').Value = '운동화'
$ws.Range('J138').Value = '25'
$ws.Range('B139').Value = '37'
$ws.Range('C139').Value = 'soom._.inc'
$ws.Range('D139').Value = '1,7'
$ws.Range('E139').Value = '바람막이'
$ws.Range('H139').Value = '면바지'
$ws.Range('I139').Value = '운동화'
$ws.Range('J139').Value = '12'
$ws.Range('B140').Value = '38'
$ws.Range('C140').Value = 'soom._.inc'
$ws.Range('D140').Value = '1,7'
$ws.Range('F140').Value = '긴팔'
$ws.Range('H140').Value = '면바지'
$ws.Range('I140').Value = '단화'
$ws.Range('J140').Value = '19'
$ws.Range('B141').Value = '39'
$ws.Range('C141').Value = 'soom._.inc'
$ws.Range('D141').Value = '1,7'
$ws.Range('F141').Value = '긴팔'
$ws.Range('H141').Value = '면바지'
$ws.Range('I141').Value = '운동화'
$ws.Range('J141').Value = '19'
$ws.Range('B142').Value = '40'
$ws.Range('C142').Value = 'soom._.inc'
$ws.Range('D142').Value = '1,7'
$ws.Range('F142').Value = '긴팔'
$ws.Range('H142').Value = '면바지'
$ws.Range('I142').Value = '운동화'
$ws.Range('J142').Value = '19'
$ws.Range('B143').Value = '41'
$ws.Range('C143').Value = 'soom._.inc'
$ws.Range('D143').Value = '1,7'
$ws.Range('F143').Value = '반팔'
$ws.Range('H143').Value = '면바지'
$ws.Range('I143').Value = '운동화'
$ws.Range('J143').Value = '28'
$ws.Range('B144').Value = '42'
$ws.Range('C144').Value = 'soom._.inc'
$ws.Range('D144').Value = '1,7'
$ws.Range('F144').Value = '반팔'
$ws.Range('G144').Value = '긴팔'
$ws.Range('H144').Value = '면바지'
$ws.Range('I144').Value = '운동화'
$ws.Range('J144').Value = '22'
$ws.Range('B145').Value = '43'
$ws.Range('C145').Value = 'soom._.inc'
$ws.Range('D145').Value = '4'
$ws.Range('E145').Value = '블레이저'
$ws.Range('F145').Value = '반팔'
$ws.Range('H145').Value = '면바지'
$ws.Range('I145').Value = '구두'
$ws.Range('J145').Value = '16'
$ws.Range('B146').Value = '44'
$ws.Range('C146').Value = 'soom._.inc'
$ws.Range('D146').Value = '1,7'
$ws.Range('F146').Value = '긴팔'
$ws.Range('H146').Value = '면바지'
$ws.Range('I146').Value = '컨버스'
$ws.Range('J146').Value = '19'
$ws.Range('B147').Value = '45'
$ws.Range('C147').Value = 'soom._.inc'
$ws.Range('D147').Value = '1,7'
$ws.Range('F147').Value = '긴팔'
$ws.Range('H147').Value = '면바지'
$ws.Range('I147').Value = '운동화'
$ws.Range('J147').Value = '19'
$ws.Range('B148').Value = '46'
$ws.Range('C148').Value = 'soom._.inc'
$ws.Range('D148').Value = '1'
$ws.Range('E148').Value = '니트집업'
$ws.Range('H148').Value = '면바지'
$ws.Range('I148').Value = '운동화'
$ws.Range('J148').Value = '19'
$ws.Range('B149').Value = '47'
$ws.Range('C149').Value = 'soom._.inc'
$ws.Range('D149').Value = '1,7'
$ws.Range('F149').Value = '긴팔'
$ws.Range('H149').Value = '면바지'
$ws.Range('I149').Value = '운동화'
$ws.Range('J149').Value = '19'
$ws.Range('B150').Value = '48'
$ws.Range('C150').Value = 'soom._.inc'
$ws.Range('D150').Value = '1,7'
$ws.Range('F150').Value = '반팔'
$ws.Range('H150').Value = '면바지'
$ws.Range('I150').Value = '컨버스'
$ws.Range('J150').Value = '28'
$ws.Range('B151').Value = '49'
$ws.Range('C151').Value = 'soom._.inc'
$ws.Range('D151').Value = '1,7'
$ws.Range('F151').Value = '셔츠'
$ws.Range('G151').Value = '반팔'
$ws.Range('H151').Value = '면바지'
$ws.Range('I151').Value = '컨버스'
$ws.Range('J151').Value = '25'
$ws.Range('B152').Value = '50'
$ws.Range('C152').Value = 'soom._.inc'
$ws.Range('D152').Value = '1'
$ws.Range('F152').Value = '반팔'
$ws.Range('H152').Value = '반바지'
$ws.Range('I152').Value = '컨버스'
$ws.Range('J152').Value = '28'
$ws.Range('B153').Value = '51'
$ws.Range('C153').Value = 'soom._.inc'
$ws.Range('D153').Value = '4'
$ws.Range('F153').Value = '반팔'
$ws.Range('H153').Value = '치마바지'
$ws.Range('I153').Value = '구두'
$ws.Range('J153').Value = '25'
$ws.Range('B154').Value = '52'
$ws.Range('C154').Value = 'soom._.inc'
$ws.Range('D154').Value = '1,7'
$ws.Range('F154').Value = '반팔'
$ws.Range('H154').Value = '면바지'
$ws.Range('I154').Value = '컨버스'
$ws.Range('J154').Value = '28'
$ws.Range('B155').Value = '53'
$ws.Range('C155').Value = 'soom._.inc'
$ws.Range('D155').Value = '1,7'
$ws.Range('F155').Value = '반팔'
$ws.Range('H155').Value = '면바지'
$ws.Range('I155').Value = '운동화'
$ws.Range('J155').Value = '28'
$ws.Range('B156').Value = '54'
$ws.Range('C156').Value = 'soom._.inc'
$ws.Range('D156').Value = '1,7'
$ws.Range('F156').Value = '반팔'
$ws.Range('H156').Value = '반바지'
$ws.Range('I156').Value = '단화'
$ws.Range('J156').Value = '28'
$ws.Range('B157').Value = '55'
$ws.Range('C157').Value = 'soom._.inc'
$ws.Range('D157').Value = '1,7'
$ws.Range('F157').Value = '긴팔'
$ws.Range('H157').Value = '면바지'
$ws.Range('I157').Value = '운동화'
$ws.Range('J157').Value = '25'
$ws.Range('B158').Value = '56'
$ws.Range('C158').Value = 'soom._.inc'
$ws.Range('D158').Value = '7'
$ws.Range('F158').Value = '반팔'
$ws.Range('H158').Value = '치마'
$ws.Range('I158').Value = '단화'
$ws.Range('J158').Value = '28'

# Reflect the final view/selection state used while entering this data
# (scrolled so row 111 is at the top, with L151 as the active cell)
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 111
    $win.ScrollColumn = 1
} catch {
    # view-state scrolling isn't critical to the data itself; ignore if unsupported
}
$ws.Range("L151").Select()
